$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section: "Column descriptions for CO2_log_exp1.csv" metadata block (rows 15-21)

# Row 15: section header (bold, matches the style used for A1/A5 headers)
$ws.Range("A15").Value = 'Column descriptions for "CO2_log_exp1.csv"'
$ws.Range("A15").Font.Bold = $true

# Row 16: time
$ws.Range("A16").Value = "time"
$ws.Range("B16").Value = "Time stamp for CO2 reading in the format YYMMDD_HHMMSS"

# Row 17: co2_a
$ws.Range("A17").Value = "co2_a"
$ws.Range("B17").Value = "CO2 ppm reading from Licor 8500 in ambient CO2 greenhouse"

# Row 18: co2_e
$ws.Range("A18").Value = "co2_e"
$ws.Range("B18").Value = "CO2 ppm reading from Licor 8500 in elevated CO2 greenhouse"

# Row 19: temp_a
$ws.Range("A19").Value = "temp_a"
$ws.Range("B19").Value = "Temperature reading from Licor 8500 in ambient CO2 greenhouse"

# Row 20: temp_e
$ws.Range("A20").Value = "temp_e"
$ws.Range("B20").Value = "Temperature reading from Licor 8500 in elevated CO2 greenhouse"

# Row 21: datenum
$ws.Range("A21").Value = "datenum"
$ws.Range("B21").Value = 'Lubridate value for time stamp in "time" column used for plotting in R '

# Update selection to match the final cursor position left by the author
$ws.Range("A21").Select()
